$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.005") into actual numbers, which would
# change their value/precision. Column D already contains plain text values
# in the source workbook (t="inlineStr"), so this keeps the same semantics.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.014.23'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.645.03'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '215.62'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").Value = '0.5091'
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").Value = '0.2572'
$ws.Range("E8").Value = '  +0.31%  '

$ws.Range("D9").Value = '0.06389'
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("D10").Value = '19.59'
$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("D11").Value = '0.07780'
$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.659.55'
$ws.Range("E12").Value = '  +1.13%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.297'
$ws.Range("E13").Value = '  +1.16%  '

$ws.Range("D14").Value = '0.5473'
$ws.Range("E14").Value = '  +0.44%  '

$ws.Range("D15").Value = '0.0₅7850'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '64.62'
$ws.Range("E16").Value = '  +1.05%  '

$ws.Range("D17").Value = '26.100.47'
$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").Value = '199.05'
$ws.Range("E19").Value = '  -1.59%  '

$ws.Range("D20").Value = '4.477'
$ws.Range("E20").Value = '  +2.29%  '

$ws.Range("D21").Value = '9.990'
$ws.Range("E21").Value = '  +1.23%  '

$ws.Range("D22").Value = '6.070'
$ws.Range("E22").Value = '  +1.86%  '

$ws.Range("D23").Value = '1.006'
$ws.Range("E23").Value = '  +0.34%  '

$ws.Range("D24").Value = '1.885'
$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").Value = '141.73'
$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("D26").Value = '0.1174'
$ws.Range("E26").Value = '  +3.77%  '

$ws.Range("D27").Value = '6.911'
$ws.Range("E27").Value = '  +2.53%  '

$ws.Range("D28").Value = '15.76'
$ws.Range("E28").Value = '  +0.66%  '

$ws.Range("D29").Value = '1.242'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = '0.05027'
$ws.Range("E30").Value = '  +1.23%  '

$ws.Range("D31").Value = '3.267'
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").Value = '3.201'
$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").Value = '1.546'
$ws.Range("E33").Value = '  +0.35%  '

$ws.Range("D34").Value = '2.365'
$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").Value = '0.9023'
$ws.Range("E35").Value = '  +1.32%  '

$ws.Range("D36").Value = '2.589'
$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("D37").Value = '1.132.11'
$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("D38").Value = '0.5496'
$ws.Range("E38").Value = '  -1.86%  '

$ws.Range("D39").Value = '0.01561'
$ws.Range("E39").Value = '  -0.21%  '

$ws.Range("D40").Value = '1.005'
$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("D41").Value = '2.544'
$ws.Range("E41").Value = '  -1.23%  '

$ws.Range("B42").Value = 'BabyDogeCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D42").Value = '0.0₈129'
$ws.Range("E42").Value = '  +10.06%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.629'
$ws.Range("E43").Value = '  -0.53%  '

$ws.Range("D44").Value = '0.8215'
$ws.Range("E44").Value = '  +2.11%  '

$ws.Range("D45").Value = '100.19'
$ws.Range("E45").Value = '  +0.60%  '

$ws.Range("D46").Value = '1.785.09'
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("D47").Value = '0.4546'
$ws.Range("E47").Value = '  +0.46%  '

$ws.Range("D48").Value = '1.008'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").Value = '54.97'
$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("D50").Value = '0.05076'
$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("D51").Value = '1.009'
$ws.Range("E51").Value = '  +0.76%  '
